$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = "14:10-14:15"
$ws.Range("C9").Value = "14:15-14:20"

$ws.Range("C16").Select()
